$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C (2nd header "Förändrad") holds a date serial that was bumped by
# one day (2023-09-12 -> 2023-09-13, i.e. 45181 -> 45182) for every data
# row (rows 2 through 329).
$ws.Range("C2:C329").Value = 45182
